$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.5
$ws.Range("F3").Value = 2.28
$ws.Range("G3").Value = 2.6
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 3.2
$ws.Range("L3").Value = 1.42
$ws.Range("Q3").Value = 1.98
$ws.Range("V3").Value = 1.37
$ws.Range("W3").Value = 1.63
$ws.Range("Z3").Value = 25
$ws.Range("AA3").Value = 65
$ws.Range("AD3").Value = 15.5
$ws.Range("I6").Value = 2.46
$ws.Range("J6").Value = 3.3
$ws.Range("K6").Value = 3.5
$ws.Range("Q6").Value = 2.2
$ws.Range("V6").Value = 1.68
$ws.Range("I7").Value = 6.2
$ws.Range("V7").Value = 1.2
$ws.Range("L8").Value = 1.33
$ws.Range("G10").Value = 1.84
$ws.Range("H10").Value = 5.5
$ws.Range("J10").Value = 3.65
$ws.Range("M10").Value = 1.07
$ws.Range("Q10").Value = 2.04
$ws.Range("W10").Value = 2.18
$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.86
$ws.Range("W11").Value = 1.57
$ws.Range("P12").Value = 1.66
$ws.Range("Q12").Value = 2.42
$ws.Range("T12").Value = 2.06
$ws.Range("G13").Value = 2.42
$ws.Range("I13").Value = 3.35
$ws.Range("K13").Value = 3.95
$ws.Range("L13").Value = 1.34
$ws.Range("N13").Value = 4.5
$ws.Range("P13").Value = 2.2
$ws.Range("W13").Value = 1.7
$ws.Range("G14").Value = 2.58
$ws.Range("K14").Value = 3.7
$ws.Range("Q14").Value = 1.99
$ws.Range("AF14").Value = 980
$ws.Range("F15").Value = 3
$ws.Range("J15").Value = 3.25
$ws.Range("H16").Value = 2.78
$ws.Range("K16").Value = 4.2
$ws.Range("S16").Value = 2.58
$ws.Range("U16").Value = 2.46
$ws.Range("G17").Value = 2.16
$ws.Range("H17").Value = 3.5
$ws.Range("K17").Value = 4.8
$ws.Range("T17").Value = 1.61
$ws.Range("U17").Value = 2.4
$ws.Range("W17").Value = 1.86
$ws.Range("G18").Value = 1.99
$ws.Range("H18").Value = 3.95
$ws.Range("Q18").Value = 1.59
$ws.Range("W18").Value = 2
$ws.Range("AN18").Value = 9.4
$ws.Range("Q19").Value = 1.56
$ws.Range("U19").Value = 2.3
$ws.Range("H20").Value = 2.92
$ws.Range("K20").Value = 3.15
$ws.Range("P20").Value = 1.56
$ws.Range("Z20").Value = 23
$ws.Range("AA20").Value = 60
$ws.Range("G22").Value = 4.6
$ws.Range("I22").Value = 2.2
$ws.Range("K22").Value = 3.55
$ws.Range("V22").Value = 1.83
$ws.Range("AB22").Value = 13.5
$ws.Range("AH22").Value = 980
$ws.Range("F23").Value = 2.18
$ws.Range("G23").Value = 2.2
$ws.Range("J23").Value = 3.25
$ws.Range("K23").Value = 3.3
$ws.Range("P23").Value = 1.68
$ws.Range("Q23").Value = 2.38
$ws.Range("W23").Value = 1.83
$ws.Range("Z23").Value = 27
$ws.Range("F24").Value = 2.22
$ws.Range("G24").Value = 2.46
$ws.Range("H24").Value = 4.1
$ws.Range("I24").Value = 4.7
$ws.Range("V24").Value = 1.27
$ws.Range("W24").Value = 1.68
$ws.Range("F25").Value = 1.81
$ws.Range("G25").Value = 1.82
$ws.Range("N25").Value = 4.6
$ws.Range("P25").Value = 2.2
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 1.79
$ws.Range("U25").Value = 2.22
$ws.Range("AH25").Value = 18
$ws.Range("AK25").Value = 17.5
$ws.Range("H26").Value = 4.2
$ws.Range("N26").Value = 3.25
$ws.Range("G27").Value = 1.73
$ws.Range("J27").Value = 3.65
$ws.Range("L27").Value = 1.47
$ws.Range("O27").Value = 1.39
$ws.Range("W27").Value = 2.36
$ws.Range("Q28").Value = 2.3
$ws.Range("V28").Value = 1.49
$ws.Range("G29").Value = 2.42
$ws.Range("I29").Value = 4.3
$ws.Range("J29").Value = 3.35
$ws.Range("Q29").Value = 1.96
$ws.Range("S29").Value = 3.25
$ws.Range("W29").Value = 1.71
$ws.Range("V30").Value = 1.34
$ws.Range("AJ30").Value = 38
$ws.Range("Q31").Value = 2.62
